# spring 24 week 8 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$newData = @(
    @(3, 5, 2, 15),
    @(8, 14, 5, 6),
    @(4, 8, 3, 12),
    @(1, 7, 2, 13),
    @(6, 4, 8, 16),
    @(5, 8, 6, 12),
    @(5, 7, 2, 13),
    @(3, 19, 4, 1),
    @(5, 7, 7, 13),
    @(1, 14, 2, 6),
    @(6, 6, 4, 14),
    @(2, 16, 3, 4),
    @(6, 8, 5, 12),
    @(4, 15, 5, 5),
    @(4, 8, 3, 12),
    @(5, 4, 6, 16),
    @(5, 3, 6, 17),
    @(6, 5, 8, 15),
    @(4, 15, 5, 5),
    @(2, 16, 1, 4),
    @(3, 6, 2, 14),
    @(4, 6, 7, 14),
    @(4, 5, 5, 15),
    @(3, 2, 4, 18),
    @(3, 14, 5, 6),
    @(6, 12, 5, 8),
    @(4, 8, 2, 12),
    @(6, 6, 5, 14),
    @(4, 5, 3, 15),
    @(2, 8, 3, 12),
    @(6, 3, 5, 17),
    @(4, 14, 7, 6),
    @(6, 6, 4, 14),
    @(6, 8, 4, 12),
    @(4, 16, 5, 4),
    @(2, 19, 4, 1)
)

$startRow = 1472
$endRow = $startRow + $newData.Count - 1

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $newData[$i]
    $ws.Range("A$r").Value = $rowVals[0]
    $ws.Range("B$r").Value = $rowVals[1]
    $ws.Range("C$r").Value = $rowVals[2]
    $ws.Range("D$r").Value = $rowVals[3]
}

# Update the view: scroll and selection to mirror the new end of data
$ws.Activate()
$newSelCell = "A$($endRow + 1)"
$ws.Range($newSelCell).Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 1484
$win.ScrollColumn = 1
